$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare new row 8 by copying the style of row 7 column A (so A8 gets s="1" like the rest)
$ws.Range("A7").Copy($ws.Range("A8"))

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2023
$ws.Range("D2").Value = 226.5
$ws.Range("F2").Value = 113.5735294117647
$ws.Range("G2").Value = 3.5
$ws.Range("H2").Value = "Detroit"
$ws.Range("I2").Value = "Indiana"
$ws.Range("J2").Value = 0.5149253731343284
$ws.Range("K2").Value = 99.55073529411763
$ws.Range("L2").Value = 113.4705882352941
$ws.Range("M2").Value = 118.2963235294118
$ws.Range("N2").Value = 73.59264705882353
$ws.Range("O2").Value = 0.4004999999999999
$ws.Range("P2").Value = 0.5714411764705882
$ws.Range("Q2").Value = 0.284375
$ws.Range("R2").Value = 12.47573529411765
$ws.Range("S2").Value = 12.27941176470589
$ws.Range("T2").Value = 0.2275845588235293
$ws.Range("U2").Value = 0.9919085538145389
$ws.Range("V2").Value = 1.074257995036756
$ws.Range("W2").Value = 11.7836500394311
$ws.Range("X2").Value = 0.3382352941176471
$ws.Range("Y2").Value = 27
$ws.Range("Z2").Value = 75.05000000000001
$ws.Range("AA2").Value = 0.3462732919254659
$ws.Range("AB2").Value = 0.4923062650374564
$ws.Range("AC2").Value = -5.305
$ws.Range("AD2").Value = 0.5731698060805751
$ws.Range("AE2").Value = 0.4497074142156863
$ws.Range("AF2").Value = 0.4960614951192897

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2023
$ws.Range("D3").Value = 244.5
$ws.Range("F3").Value = 116.3455882352941
$ws.Range("G3").Value = 5.5
$ws.Range("H3").Value = "Atlanta"
$ws.Range("I3").Value = "Minnesota"
$ws.Range("J3").Value = 0.5220588235294117
$ws.Range("K3").Value = 100.2169117647059
$ws.Range("L3").Value = 115.4102941176471
$ws.Range("M3").Value = 115.5125
$ws.Range("N3").Value = 75.175
$ws.Range("O3").Value = 0.3585514705882353
$ws.Range("P3").Value = 0.5865661764705885
$ws.Range("Q3").Value = 0.2554705882352941
$ws.Range("R3").Value = 12.09044117647059
$ws.Range("S3").Value = 12.55147058823529
$ws.Range("T3").Value = 0.2114779411764706
$ws.Range("U3").Value = 1.016118674544053
$ws.Range("V3").Value = 1.025677006274918
$ws.Range("W3").Value = 10.88176881783746
$ws.Range("X3").Value = 0.5
$ws.Range("Y3").Value = 48
$ws.Range("Z3").Value = 76.25
$ws.Range("AA3").Value = 0.525
$ws.Range("AB3").Value = 0.4997640803003665
$ws.Range("AC3").Value = 2.04
$ws.Range("AD3").Value = 0.07389385433594325
$ws.Range("AE3").Value = 0.05870506535947714
$ws.Range("AF3").Value = 0.4913680646396232

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2023
$ws.Range("D4").Value = 224.5
$ws.Range("F4").Value = 114.8424688057041
$ws.Range("G4").Value = 1.5
$ws.Range("H4").Value = "Dallas"
$ws.Range("I4").Value = "Memphis"
$ws.Range("J4").Value = 0.4948694029850746
$ws.Range("K4").Value = 98.06198752228164
$ws.Range("L4").Value = 116.329055258467
$ws.Range("M4").Value = 114.0671122994653
$ws.Range("N4").Value = 76.18863636363636
$ws.Range("O4").Value = 0.4255871212121212
$ws.Range("P4").Value = 0.5814988859180037
$ws.Range("Q4").Value = 0.2928224153297683
$ws.Range("R4").Value = 11.06726827094474
$ws.Range("S4").Value = 12.24329322638146
$ws.Range("T4").Value = 0.2195296345811051
$ws.Range("U4").Value = 1.00299099393628
$ws.Range("V4").Value = 0.9867825902036287
$ws.Range("W4").Value = 11.66801595401862
$ws.Range("X4").Value = 0.553030303030303
$ws.Range("Y4").Value = 49
$ws.Range("Z4").Value = 75.25
$ws.Range("AA4").Value = 0.5091954022988505
$ws.Range("AB4").Value = 0.4875867848749801
$ws.Range("AC4").Value = 4.245
$ws.Range("AD4").Value = 0.2855650450752738
$ws.Range("AE4").Value = 0.1673476647603486
$ws.Range("AF4").Value = 0.4942606702759253

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2023
$ws.Range("D5").Value = 225.5
$ws.Range("F5").Value = 112.9149616368287
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = "Miami"
$ws.Range("I5").Value = "Utah"
$ws.Range("J5").Value = 0.5121870882740448
$ws.Range("K5").Value = 97.77088661551576
$ws.Range("L5").Value = 114.7586423699915
$ws.Range("M5").Value = 115.2360187553282
$ws.Range("N5").Value = 76.54703751065642
$ws.Range("O5").Value = 0.4176441815856777
$ws.Range("P5").Value = 0.5775859974424553
$ws.Range("Q5").Value = 0.2735002131287297
$ws.Range("R5").Value = 12.48536871270247
$ws.Range("S5").Value = 12.53610400682012
$ws.Range("T5").Value = 0.2136608589087809
$ws.Range("U5").Value = 0.9861568701906432
$ws.Range("V5").Value = 1.031404075660796
$ws.Range("W5").Value = 9.655854799659785
$ws.Range("X5").Value = 0.5035166240409207
$ws.Range("Y5").Value = 36.5
$ws.Range("Z5").Value = 75.25
$ws.Range("AA5").Value = 0.5417620137299771
$ws.Range("AB5").Value = 0.5022842464371763
$ws.Range("AC5").Value = 4.95
$ws.Range("AD5").Value = 0.1596861005467862
$ws.Range("AE5").Value = 0.1226189746732026
$ws.Range("AF5").Value = 0.5100891362589994

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 2023
$ws.Range("D6").Value = 231.5
$ws.Range("F6").Value = 114.0607989464443
$ws.Range("G6").Value = 12.5
$ws.Range("H6").Value = "Houston"
$ws.Range("I6").Value = "Boston"
$ws.Range("J6").Value = 0.5151515151515151
$ws.Range("K6").Value = 98.19494073748899
$ws.Range("L6").Value = 114.9314091308165
$ws.Range("M6").Value = 116.2337686567164
$ws.Range("N6").Value = 77.63278094820018
$ws.Range("O6").Value = 0.4220735294117646
$ws.Range("P6").Value = 0.5775081211589114
$ws.Range("Q6").Value = 0.2687833625987709
$ws.Range("R6").Value = 12.66546312554873
$ws.Range("S6").Value = 11.09462247585601
$ws.Range("T6").Value = 0.2069004609306409
$ws.Range("U6").Value = 0.9961641829383777
$ws.Range("V6").Value = 1.015156363027759
$ws.Range("W6").Value = 11.62927082437351
$ws.Range("X6").Value = 0.4575285338015803
$ws.Range("Y6").Value = 39
$ws.Range("Z6").Value = 76.15
$ws.Range("AA6").Value = 0.4712918660287082
$ws.Range("AB6").Value = 0.5036654602801898
$ws.Range("AC6").Value = -0.6200000000000001
$ws.Range("AD6").Value = 0.03328213792044232
$ws.Range("AE6").Value = 0.03804381127450981
$ws.Range("AF6").Value = 0.5234898583329538

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 2023
$ws.Range("D7").Value = 236
$ws.Range("F7").Value = 115.7381474978051
$ws.Range("G7").Value = 4.5
$ws.Range("H7").Value = "GoldenState"
$ws.Range("I7").Value = "Phoenix"
$ws.Range("J7").Value = 0.5378787878787878
$ws.Range("K7").Value = 99.21947980684811
$ws.Range("L7").Value = 115.8508340649693
$ws.Range("M7").Value = 114.5007462686567
$ws.Range("N7").Value = 76.36245610184375
$ws.Range("O7").Value = 0.4246378402107111
$ws.Range("P7").Value = 0.5852208077260754
$ws.Range("Q7").Value = 0.2337589991220369
$ws.Range("R7").Value = 12.74511633011413
$ws.Range("S7").Value = 12.2415825285338
$ws.Range("T7").Value = 0.2076905728709394
$ws.Range("U7").Value = 1.010813515264673
$ws.Range("V7").Value = 1.072352460686536
$ws.Range("W7").Value = 11.55745725821847
$ws.Range("X7").Value = 0.5334723441615452
$ws.Range("Y7").Value = 52.5
$ws.Range("Z7").Value = 76.65
$ws.Range("AA7").Value = 0.4495798319327731
$ws.Range("AB7").Value = 0.4948403696302915
$ws.Range("AC7").Value = 6.23
$ws.Range("AD7").Value = 0.08958064490401504
$ws.Range("AE7").Value = 0.07502117310062878
$ws.Range("AF7").Value = 0.5141165355800927

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 2023
$ws.Range("D8").Value = 244.5
$ws.Range("F8").Value = 118.5463591135233
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = "Sacramento"
$ws.Range("I8").Value = "Milwaukee"
$ws.Range("J8").Value = 0.5001148105625717
$ws.Range("K8").Value = 99.93351424694708
$ws.Range("L8").Value = 117.6110244233379
$ws.Range("M8").Value = 114.3529963817278
$ws.Range("N8").Value = 77.68224785165086
$ws.Range("O8").Value = 0.4320945273631841
$ws.Range("P8").Value = 0.5967804161013115
$ws.Range("Q8").Value = 0.272714156490276
$ws.Range("R8").Value = 12.22062415196744
$ws.Range("S8").Value = 11.25609452736319
$ws.Range("T8").Value = 0.2013562302125735
$ws.Range("U8").Value = 1.035339380904134
$ws.Range("V8").Value = 1.042155439635926
$ws.Range("W8").Value = 13.08270057062465
$ws.Range("X8").Value = 0.6612392582541836
$ws.Range("Y8").Value = 44
$ws.Range("Z8").Value = 76
$ws.Range("AA8").Value = 0.6287625418060201
$ws.Range("AB8").Value = 0.4831006368827571
$ws.Range("AC8").Value = -1.02
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0.4918769531364166
